# "Add Impediment Backlog Change Team Wiki"
#
# 1) Rename Sheet1 -> "Project Backlog"
# 2) Replace the old empty Sheet2 / Sheet3 with a populated "Impediment Backlog"
#    sheet (second tab) and make it the active tab.
# 3) Widen the Project Backlog title bar from A1:F1 to A1:G1.
# 4) Fill in the new Impediment Backlog content/formatting.

$wb = $excel.ActiveWorkbook

# --- Sheets: add the replacement sheet, drop the old empties, reposition ----
$newSheet = $wb.Worksheets.Add()
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))

$wb.Worksheets.Item("Sheet1").Name = "Project Backlog"
$wb.Worksheets.Item("Sheet4").Name = "Impediment Backlog"

$ws1 = $wb.Worksheets.Item("Project Backlog")
$ws2 = $wb.Worksheets.Item("Impediment Backlog")

# --- Project Backlog: widen the title row from A1:F1 to A1:G1 --------------
$ws1.Range("A1:F1").UnMerge()
$ws1.Range("A1:G1").Merge()
$ws1.Range("A1:G1").Select()

# --- Impediment Backlog: content --------------------------------------------
# Write in this order so newly-introduced shared strings land in the same
# sequence as the target workbook (Status, Group 2 Impediment, the two
# impediment descriptions, "Not solve", then the SVN-convention note).
$ws2.Range("B2").Value = "Status"
$ws2.Range("A1").Value = "Group 2 Impediment"
$ws2.Range("C3").Value = "Some daily meeting late more than 20 minutes"
$ws2.Range("B3").Value = "Not solve"
$ws2.Range("C4").Value = "Not follow convention when commit code to SVN"

$ws2.Range("A2").Value = "No."
$ws2.Range("C2").Value = "Description"
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "Not solve"

# --- Impediment Backlog: formatting -----------------------------------------
$ws2.Range("A1:C1").Merge()
$ws2.Range("A1:C1").HorizontalAlignment = -4108

$ws2.Range("B3").Font.Color = 255
$ws2.Range("B4").Font.Color = 255

$ws2.Columns.Item(2).ColumnWidth = 10.307291666666666
$ws2.Columns.Item(3).ColumnWidth = 44.7

$ws2.PageSetup.Orientation = 1

$ws2.Range("C5").Select()

# --- Make "Impediment Backlog" the active tab --------------------------------
$ws2.Activate()
